$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "What is the difference between window, screen, and document in Java" +
#    "S" + "cript?"  ->  single run
#    "What is the difference between window, screen, and document in JavaScript?"
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
$full12 = $p12.Range
$textEnd12 = $full12.End - 1
$rAll12 = $d.Range($full12.Start, $textEnd12)
# Force a real text change so the engine rewrites the run (avoids a no-op
# when old/new text happen to render identically) and ends up as one run.
$rAll12.Text = "TEMP_PLACEHOLDER_12"
$p12b = $d.Paragraphs(12)
$full12b = $p12b.Range
$textEnd12b = $full12b.End - 1
$rAll12b = $d.Range($full12b.Start, $textEnd12b)
$rAll12b.Text = "What is the difference between window, screen, and document in JavaScript?"

# ---------------------------------------------------------------------------
# 2) "Write a blog about objects and its internal representation in Java" +
#    "Script"  ->
#    "Write a blog about objects and its internal representation in JavaScript" + "?"
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$full10 = $p10.Range
$textEnd10 = $full10.End - 1
$rAll10 = $d.Range($full10.Start, $textEnd10)
$rAll10.Text = "TEMP_PLACEHOLDER_10"
$p10b = $d.Paragraphs(10)
$full10b = $p10b.Range
$textEnd10b = $full10b.End - 1
$rAll10b = $d.Range($full10b.Start, $textEnd10b)
$rAll10b.Text = "Write a blog about objects and its internal representation in JavaScript?"

# ---------------------------------------------------------------------------
# 3) Insert a new paragraph right after the "Write a blog ... ?" paragraph
#    that contains a hyperlink to the JavaScript objects blog post.
# ---------------------------------------------------------------------------
$p10c = $d.Paragraphs(10)
$p10c.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs(11)
$rNew1 = $pNew1.Range
$rNew1.Text = "PLACEHOLDER_LINK_1"
$pNew1b = $d.Paragraphs(11)
$rSel1 = $pNew1b.Range
$rSel1.MoveEnd(1, -1)
$h1 = $d.Hyperlinks.Add($rSel1, "https://rajaking-alm.medium.com/objects-and-its-internal-representation-in-javascript-8785f57acd94", "", "", "https://rajaking-alm.medium.com/objects-and-its-internal-representation-in-javascript-8785f57acd94")
$h1.Range.Font.Name = "Arial"

# ---------------------------------------------------------------------------
# 4) Insert a new paragraph right after the "What is the difference ...?"
#    paragraph that contains a hyperlink to the window/screen/document post.
#    (It is now at index 13, since step 3 inserted a new paragraph before it.)
# ---------------------------------------------------------------------------
$p13 = $d.Paragraphs(13)
$p13.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs(14)
$rNew2 = $pNew2.Range
$rNew2.Text = "PLACEHOLDER_LINK_2"
$pNew2b = $d.Paragraphs(14)
$rSel2 = $pNew2b.Range
$rSel2.MoveEnd(1, -1)
$h2 = $d.Hyperlinks.Add($rSel2, "https://rajaking-alm.medium.com/difference-between-window-screen-and-document-in-javascript-44b74193edbf", "", "", "https://rajaking-alm.medium.com/difference-between-window-screen-and-document-in-javascript-44b74193edbf")
$h2.Range.Font.Name = "Arial"

Write-Host "Paragraphs count:" $d.Paragraphs.Count
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    Write-Host $i ": [" $p.Range.Text "]"
}
